# Bevill State Community College Organizations workbook update
# - Swap "Organization Name" (A) and "Categories" (B) columns
# - Rename several headers
# - Add new "Tiktok Link" column (M)
# - Resize several columns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 70

# --- 1) Swap the data held in columns A and B (rows 2..70) -----------------
# Column A currently holds the organization name, column B the category.
# After the edit, A should hold the category and B the organization name.
for ($r = 2; $r -le $lastRow; $r++) {
    $orgName  = $ws.Cells.Item($r, 1).Value2
    $category = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $category
    $ws.Cells.Item($r, 2).Value = $orgName
}

# --- 2) Update the header row ----------------------------------------------
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "Organization Name"
$ws.Range("C1").Value = "Organization Link"
$ws.Range("D1").Value = "Logo Link"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Email"
$ws.Range("G1").Value = "Phone Number"
$ws.Range("H1").Value = "Linkedin Link"
$ws.Range("I1").Value = "Instagram Link"
$ws.Range("J1").Value = "Facebook Link"
$ws.Range("K1").Value = "Twitter Link"
$ws.Range("L1").Value = "Youtube Link"

# --- 3) Add the new "Tiktok Link" column (M) --------------------------------
# Copy the style of an existing header cell so the new header matches
# the look (bold, centered, bordered) of the rest of row 1.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("M1").Value = "Tiktok Link"

# --- 4) Resize columns -------------------------------------------------------
# ColumnWidth values are offset by 0.83 from the stored/display character
# width used by the OOXML <col width="..."/> attribute on this runtime, so
# subtract 0.83 from each desired stored width before assigning it.
$ws.Columns.Item(1).ColumnWidth  = 18 - 0.83    # A
$ws.Columns.Item(2).ColumnWidth  = 46 - 0.83    # B
$ws.Columns.Item(3).ColumnWidth  = 50 - 0.83    # C
$ws.Columns.Item(4).ColumnWidth  = 11 - 0.83    # D
$ws.Columns.Item(5).ColumnWidth  = 13 - 0.83    # E
$ws.Columns.Item(6).ColumnWidth  = 7 - 0.83     # F
$ws.Columns.Item(7).ColumnWidth  = 14 - 0.83    # G
$ws.Columns.Item(8).ColumnWidth  = 15 - 0.83    # H
$ws.Columns.Item(9).ColumnWidth  = 16 - 0.83    # I
$ws.Columns.Item(10).ColumnWidth = 15 - 0.83    # J
$ws.Columns.Item(11).ColumnWidth = 14 - 0.83    # K
$ws.Columns.Item(12).ColumnWidth = 14 - 0.83    # L
$ws.Columns.Item(13).ColumnWidth = 13 - 0.83    # M
